# 길드전_답지.xlsx - "Add files via upload"
#
# The commit appends three trailing columns (H: 날짜/date, I: 상대 길드/opponent
# guild, J: 기준/standard) to every data row from row 475 through row 556 of
# Sheet1. All of the new cells share the same values:
#   H -> 260131 (number)
#   I -> 밤빛   (text)
#   J -> 방어   (text)
# and use the workbook's existing centered-alignment cell style (the same
# style already applied to every other populated cell on the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 475
$lastRow = 556

# --- H column: numeric date-like value -------------------------------------
$rngH = $ws.Range("H$firstRow`:H$lastRow")
$rngH.Value = 260131
$rngH.HorizontalAlignment = -4108   # xlCenter
$rngH.VerticalAlignment = -4108     # xlCenter

# --- I column: opponent guild name ------------------------------------------
$rngI = $ws.Range("I$firstRow`:I$lastRow")
$rngI.Value = "밤빛"
$rngI.HorizontalAlignment = -4108
$rngI.VerticalAlignment = -4108

# --- J column: standard/type label ------------------------------------------
$rngJ = $ws.Range("J$firstRow`:J$lastRow")
$rngJ.Value = "방어"
$rngJ.HorizontalAlignment = -4108
$rngJ.VerticalAlignment = -4108

# --- Restore the view state recorded in the saved workbook -----------------
# (best-effort; selection is what actually round-trips to the sheetView XML)
[void]$ws.Range("N552").Select()
